# Bump IG version to 1.8.2 and update publish date, plus fill in the
# missing invariant text for the root "Extension" element row (AJ1) on
# the Elements sheet, matching how the same invariant text already
# appears for the "Extension.extension" row.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$metadata.Range("B3").Value = "1.8.2"
$metadata.Range("B8").Value = "2023-09-01T14:45:29-04:00"

$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
